$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# ---------------------------------------------------------------------------
# SATURDAY section (Lassonde C mic setup) - rows 149-153 first (matches the
# order the original author typed them in, which drives shared-string
# allocation order), then the "SATURDAY" separator header row 148, then the
# cancellation note row 154.
# ---------------------------------------------------------------------------

# Row 149
$ws.Range("A149").Value = "Setup Mic"
$ws.Range("B149").Value = 42623
$ws.Range("C149").Value = "1715"
$ws.Range("D149").Value = "LAS"
$ws.Range("E149").Value = "C"
$ws.Range("F149").Value = "Take cart with mixer, 2 wireless mics and 2 mic stands from Lassonde 1011 storeroom (across from Lassonde A). Go to Lassonde C classroom (class starts at 5:30 pm but be there early in case previous class ends early). "
$ws.Rows.Item(149).RowHeight = 60

# Row 150
$ws.Range("A150").Value = "Other"
$ws.Range("B150").Value = 42623
$ws.Range("C150").Value = "1715"
$ws.Range("D150").Value = "LAS"
$ws.Range("E150").Value = "C"
$ws.Range("F150").Value = 'Log in as 5065*0 on touchscreen. (First level bar is your wireless handheld mic volume). Plug in mic cable from output of mixer to mic input on podium (XLR jack just above VHS machine in podium). Ramp up volume a bit on "Microphone 2" on touchscreen to medium volume to get level.'
$ws.Rows.Item(150).RowHeight = 75

# Row 151
$ws.Range("A151").Value = "Other"
$ws.Range("B151").Value = 42623
$ws.Range("C151").Value = "1715"
$ws.Range("D151").Value = "LAS"
$ws.Range("E151").Value = "C"
$ws.Range("F151").Value = 'Plug in power cord from cart on to power outlet on left side of podium (to left of document camera). Turn on mixer. Turn on wireless microphone receivers on cart (NOTE: DO NOT PRESS "SYNC" BUTTON" - POWER BUTTON IS FIRST BUTTON TO THE RIGHT ON RECEIVER). '
$ws.Rows.Item(151).RowHeight = 75

# Row 152
$ws.Range("A152").Value = "Other"
$ws.Range("B152").Value = 42623
$ws.Range("C152").Value = "1715"
$ws.Range("D152").Value = "LAS"
$ws.Range("E152").Value = "C"
$ws.Range("F152").Value = 'Press "MUTE" button on wireless microphones to turn on mics. Adjust volume by adjusting volume on mixer (inputs 1 and 2). Also you can adjust volume on first volume control bar on touchscreen. (NOTE: VOLUME ON TOUCHSCREEN MUST BE RAMPED UP OR DOWN INITIALLY TO GET ANY VOLUME - First volume bar).'
$ws.Rows.Item(152).RowHeight = 90

# Row 153
$ws.Range("A153").Value = "Other"
$ws.Range("B153").Value = 42623
$ws.Range("C153").Value = "1715"
$ws.Range("D153").Value = "LAS"
$ws.Range("E153").Value = "C"
$ws.Range("F153").Value = "Once volumes are set, place one mic stand with mic halfway up aisle on right and one mic stand with mic halfway up aisle on left. Demo volume controls to prof. and demo PC. Leave microphone bags with milk carton on cart in room. PLEASE FIND OUT END TIME OF CLASS FROM PROF. AND TELL MASI AS MICROPHONES ARE EXPENSIVE. TELL PROF. TO STAY WITH MICS UNTIL THEY ARE PICKED UP. TELL HIM TO CALL ext 55800   WHEN DONE (use phone in classroom)."
$ws.Rows.Item(153).RowHeight = 120

# Row 148 - separator/header row, copy format from an existing separator row
$ws.Range("A143:F143").Copy()
$ws.Range("A148:F148").PasteSpecial(-4122)
$ws.Range("B148").Value = "SATURDAY"

# Row 154 - note row (no C/D/E cells used)
$ws.Range("A154").Value = "MASI - THIS IS NOT REAL SAT"
$ws.Range("B154").Value = "FOR MY RECORDS JL"
$ws.Range("F154").Value = "LASSONDE C GOT CANCELLED LAST MINUTE BUT WANT TO SAVE INSTRUCTIONS. JEANNINE"
$ws.Rows.Item(154).RowHeight = 30

# ---------------------------------------------------------------------------
# MONDAY section - rows 159-162
# ---------------------------------------------------------------------------

# Row 159 - separator/header row
$ws.Range("A143:F143").Copy()
$ws.Range("A159:F159").PasteSpecial(-4122)
$ws.Range("B159").Value = "MONDAY"

# Row 160
$ws.Range("A160").Value = "Demo"
$ws.Range("B160").Value = 42625
$ws.Range("C160").Value = "1900"
$ws.Range("D160").Value = "CB"
$ws.Range("E160").Value = "129"

# Row 161
$ws.Range("A161").Value = "Demo"
$ws.Range("B161").Value = 42625
$ws.Range("C161").Value = "1900"
$ws.Range("D161").Value = "CLH"
$ws.Range("E161").Value = "J"

# Row 162
$ws.Range("A162").Value = "Demo"
$ws.Range("B162").Value = 42625
$ws.Range("C162").Value = "1900"
$ws.Range("D162").Value = "CLH"
$ws.Range("E162").Value = "K"

# ---------------------------------------------------------------------------
# TUESDAY section - rows 166-170
# ---------------------------------------------------------------------------

# Row 166 - separator/header row
$ws.Range("A143:F143").Copy()
$ws.Range("A166:F166").PasteSpecial(-4122)
$ws.Range("B166").Value = "TUESDAY"

# Row 167
$ws.Range("A167").Value = "Pickup PC"
$ws.Range("B167").Value = 42626
$ws.Range("C167").Value = "1830"
$ws.Range("D167").Value = "SC"
$ws.Range("E167").Value = "MDR"
$ws.Range("F167").Value = "Pick up PC and Projector carts from Stong Master's Dining Room. Pick up all cables, ac cords and matts and return to Bethune 201 storeroom. PLEASE LEAVE PORTABLE SCREEN IN ROOM. Stong MDR key is in CB 121A storeroom. STONG MDR is just to the right of Stong Dining Hall - go thru two black doors."
$ws.Rows.Item(167).RowHeight = 75

# Row 168
$ws.Range("A168").Value = "Pickup Small PA"
$ws.Range("B168").Value = 42626
$ws.Range("C168").Value = "1830"
$ws.Range("D168").Value = "SC"
$ws.Range("E168").Value = "MDR"
$ws.Range("F168").Value = "Pick up Small Speaker on cart and return to Bethune 201 storeroom."
$ws.Rows.Item(168).RowHeight = 30

# Row 169
$ws.Range("A169").Value = "Pickup Mic"
$ws.Range("B169").Value = 42626
$ws.Range("C169").Value = "1830"
$ws.Range("D169").Value = "SC"
$ws.Range("E169").Value = "MDR"
$ws.Range("F169").Value = "Pick up Lecturn mic stand, mic and cable and return to Bethune 201 storeroom."
$ws.Rows.Item(169).RowHeight = 30

# Row 170
$ws.Range("A170").Value = "Demo"
$ws.Range("B170").Value = 42626
$ws.Range("C170").Value = "1900"
$ws.Range("D170").Value = "CLH"
$ws.Range("E170").Value = "J"

# ---------------------------------------------------------------------------
# Sheet view bookkeeping to match the final selection state.
# ---------------------------------------------------------------------------
$ws.Range("A170").Select()
